# Fruta / hortaliza, semanal
# Insert two new daily price records for Granada (Vega Modelo de Temuco)
# at rows 115-116, pushing the existing data (previously at rows 115-209)
# down to rows 117-211.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 115; everything from the old row 115
# onward shifts down by two rows (old 115 -> 117, old 209 -> 211).
$ws.Rows("115:116").Insert()

# --- New row 115 ---
$ws.Range("A115").Value = 10
$ws.Range("B115").Value = "Vega Modelo de Temuco"
$ws.Range("C115").Value = "La Araucanía"
$ws.Range("D115").Value = 45040
$ws.Range("E115").Value = 9
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100104
$ws.Range("H115").Value = "Frutos de pepita"
$ws.Range("I115").Value = 100104001
$ws.Range("J115").Value = "Granada"
$ws.Range("K115").Value = "Wonderfull"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 25
$ws.Range("N115").Value = 20000
$ws.Range("O115").Value = 20000
$ws.Range("P115").Value = 20000
$ws.Range("Q115").Value = "$/bandeja 15 kilos granel"
$ws.Range("R115").Value = "Provincia de Limarí"
$ws.Range("S115").Value = 1333
$ws.Range("T115").Value = 15

# --- New row 116 ---
$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value = 45040
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100104
$ws.Range("H116").Value = "Frutos de pepita"
$ws.Range("I116").Value = 100104001
$ws.Range("J116").Value = "Granada"
$ws.Range("K116").Value = "Wonderfull"
$ws.Range("L116").Value = "Segunda"
$ws.Range("M116").Value = 30
$ws.Range("N116").Value = 16000
$ws.Range("O116").Value = 16000
$ws.Range("P116").Value = 16000
$ws.Range("Q116").Value = "$/bandeja 15 kilos granel"
$ws.Range("R116").Value = "Provincia de Limarí"
$ws.Range("S116").Value = 1067
$ws.Range("T116").Value = 15

# Make sure the D column (date) keeps the date/time number format used
# throughout the rest of the column.
$ws.Range("D115:D116").NumberFormat = $ws.Range("D117").NumberFormat
